$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column Q with a 2020 data point, mirroring the formatting
# already used by the neighbouring 2019 column (P) and by the thin
# separator row (row 3).

# Row 3: empty separator cell, same style as P3
$ws.Range("P3").Copy($ws.Range("Q3"))

# Row 4: year header 2020, same style as P4 but vertical-top aligned
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").VerticalAlignment = -4160
$ws.Range("Q4").Value = 2020

# Row 5: data value, keeps the pre-existing style (already s=13), just needs a value
$ws.Range("Q5").Value = 1.1000000000000001

# Row 6: data value with a bordered / numeric-formatted style
$ws.Range("D3").Copy($ws.Range("Q6"))
$ws.Range("Q6").NumberFormat = "0.0"
$ws.Range("Q6").Value = 7

# Restore the cursor / selection position as recorded in the saved file
[void]$ws.Range("J22").Select()
